$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The header cell A1 was edited from "technology" to "Technology"
$ws.Range("A1").Value = "Technology"

# Move the active selection, matching the saved selection in the diff
$ws.Range("G8").Select()
